# Add Area / Atotal columns to the Q (discharge) table on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 11): new "Area" / "Atotal" headers, plus a
#     compact summary pair (J11/K11) mirroring Atotal/Qtotal. ---
$ws.Range("G11").Value = "Area"
$ws.Range("H11").Value = "Atotal"
$ws.Range("J11").Value = "Atotal"
$ws.Range("K11").Value = "Qtotal"

# --- Row 12: per-segment area plus the running Atotal/Qtotal summary. ---
$ws.Range("G12").Formula = "=(D12-0)*B12/100"
$ws.Range("H12").Formula = "=SUM(G12:G21)"
$ws.Range("J12").Formula = "=H12"
$ws.Range("K12").Formula = "=F12"

# --- Rows 13-25: per-segment area formulas (mirrors the existing
#     E-column "segment Q" formulas, one row lower on the D reference). ---
for ($r = 13; $r -le 25; $r++) {
    $prev = $r - 1
    $ws.Range("G$r").Formula = "=(D$r-D$prev)*B$r/100"
}

# --- Selection moves to the new summary cells. ---
$ws.Range("J12:K12").Select()
